$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newValues = @(
    '30-26=',
    '87-39=',
    '18+58=',
    '84-68=',
    '90-15=',
    '60-37=',
    '73-8=',
    '19+24=',
    '33-26=',
    '52-5=',
    '96-59=',
    '8+63=',
    '28+48=',
    '39+17=',
    '93-28=',
    '19+62=',
    '50-34=',
    '38+18=',
    '7+67=',
    '71-17=',
    '92-87=',
    '80-38=',
    '90-48=',
    '63-9=',
    '78-59=',
    '73-44=',
    '76+15=',
    '24-7=',
    '26-19=',
    '36+55=',
    '71-59=',
    '91-14=',
    '71-33=',
    '67+8=',
    '32-13=',
    '86+6=',
    '52-47=',
    '16+77=',
    '30-2=',
    '19+57=',
    '6+45=',
    '49+45=',
    '51-17=',
    '58+34=',
    '45+9=',
    '29+53=',
    '28+16=',
    '29+16=',
    '63+18=',
    '77+5=',
    '87-78=',
    '93-66=',
    '7+36=',
    '7+14=',
    '93-24=',
    '69+5=',
    '35+8=',
    '9+7=',
    '19+77=',
    '16+79=',
    '82-8=',
    '70-21=',
    '8+4=',
    '66-49=',
    '63+9=',
    '45-7=',
    '96-57=',
    '47+8=',
    '70-64=',
    '28+53=',
    '91-55=',
    '28+8=',
    '7+55=',
    '58+17=',
    '68+15=',
    '24+68=',
    '17+79=',
    '52-6=',
    '90-58=',
    '51-43=',
    '60-51=',
    '57+19=',
    '64-36=',
    '16+55=',
    '27+67=',
    '56+18=',
    '47-29=',
    '50-7=',
    '49+42=',
    '22+29=',
    '49+45=',
    '81-47=',
    '44-9=',
    '15+79=',
    '5+66=',
    '80-44=',
    '3+89=',
    '83-47=',
    '27+5=',
    '93-6='
)

$idx = 0
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    $row = $t.Rows.Item($r)
    for ($c = 1; $c -le $row.Cells.Count; $c++) {
        $cell = $row.Cells.Item($c)
        $cell.Range.Text = $newValues[$idx]
        $idx = $idx + 1
    }
}

Write-Output ("Updated " + $idx + " cells")